$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "4T"
$ws.Range("B2").Value = 0.997

$ws.Range("A3").Value = "5A"
$ws.Range("B3").Value = 0.99075

$ws.Range("A4").Value = "5H"
$ws.Range("B4").Value = 1.0125

$ws.Range("A5").Value = "5Q"
$ws.Range("B5").Value = 1.0295

$ws.Range("A6").Value = "5S"
$ws.Range("B6").Value = 1.023666667

$ws.Range("A8").Value = "21L"
$ws.Range("B8").Value = 1.002333333

$ws.Range("A9").Value = "24I"
$ws.Range("B9").Value = 0.990666667

$ws.Range("A10").Value = "27C"
$ws.Range("B10").Value = 0.793

$ws.Range("A11").Value = "27L"
$ws.Range("B11").Value = 0.844

$ws.Range("A12").Value = "27M"
$ws.Range("B12").Value = 0.852

$ws.Range("A13").Value = "27Q"
$ws.Range("B13").Value = 0.826
